# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" column (D) for every row that was
# just (re-)handed off, on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$handoffTime_zhcn = "2016-03-09 18:23:38"
foreach ($r in $rows) {
    $ws_zhcn.Range("D" + $r).Value = $handoffTime_zhcn
}

$ws_dede = $wb.Worksheets.Item("de-de")
$handoffTime_dede = "2016-03-09 18:23:43"
foreach ($r in $rows) {
    $ws_dede.Range("D" + $r).Value = $handoffTime_dede
}
